$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.695.59'
$ws.Range("E2").Value = '  -0.74%  '

$ws.Range("D3").Value = '2.526.19'
$ws.Range("E3").Value = '  -2.16%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("E7").Value = '  -1.33%  '

$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0803'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.41%  '

$ws.Range("E13").Value = '  +0.01%  '

$ws.Range("D14").Value = '2.917.04'
$ws.Range("E14").Value = '  -1.97%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.82%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.534.57'
$ws.Range("E16").Value = '  -2.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.808'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.59%  '

$ws.Range("D18").Value = '42.677.91'
$ws.Range("E18").Value = '  -0.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.53%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0951'
$ws.Range("E20").Value = '  -1.99%  '

$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.67'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.63%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.51'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '157.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.73'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.18%  '

$ws.Range("E33").Value = '  +10.39%  '

$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.64'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.30%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0782'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.96%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.59%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.98'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.54%  '

$ws.Range("E39").Value = '  -1.63%  '

$ws.Range("E40").Value = '  -1.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.95%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.98%  '

$ws.Range("E43").Value = '  +0.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0300'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.67%  '

$ws.Range("E45").Value = '  +0.09%  '

$ws.Range("D46").Value = '2.002.60'
$ws.Range("E46").Value = '  -0.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.75%  '

$ws.Range("D48").Value = '2.770.78'
$ws.Range("E48").Value = '  -1.97%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.189'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '78.95'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.25%  '
